$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new column (evaluator_partial_correctness) header to sheet1 FIRST, so
#     that copying the sheet below also carries the style/header onto the new sheets ---
$ws1.Range("D1").Copy() | Out-Null
$ws1.Range("E1").PasteSpecial(-4122) | Out-Null
$ws1.Range("E1").Value = "evaluator_partial_correctness"

# --- Create the two new sheets by copying sheet1 (keeps the header formatting/style
#     without the engine stamping a fresh "baseColWidth" default onto a blank new sheet) ---
$ws1.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$ws2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2.Name = "o_20"
$ws1.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "o_20_jumbled"

# --- Text blobs (prompts / solutions / llm responses) ---
$prompt1 = @"
 Given is the adjacency matrix for a weighted undirected graph containing 16 nodes labelled A to P. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   



what is the least cost path from node A to node P?

   A B C D E F G H I J K L M N O P
 A 0 2 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 B 2 0 5 0 0 1 0 0 0 0 0 0 0 0 0 0
 C 0 5 0 3 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 3 0 0 0 0 2 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 5 0 0 5 0 0 0 0 0 0 0
 F 0 1 0 0 5 0 4 0 0 4 0 0 0 0 0 0
 G 0 0 0 0 0 4 0 0 0 0 2 0 0 0 0 0
 H 0 0 0 2 0 0 0 0 0 0 0 4 0 0 0 0
 I 0 0 0 0 5 0 0 0 0 0 0 0 3 0 0 0
 J 0 0 0 0 0 4 0 0 0 0 1 0 0 1 0 0
 K 0 0 0 0 0 0 2 0 0 1 0 3 0 0 2 0
 L 0 0 0 0 0 0 0 4 0 0 3 0 0 0 0 3
 M 0 0 0 0 0 0 0 0 3 0 0 0 0 1 0 0
 N 0 0 0 0 0 0 0 0 0 1 0 0 1 0 2 0
 O 0 0 0 0 0 0 0 0 0 0 2 0 0 2 0 5
 P 0 0 0 0 0 0 0 0 0 0 0 3 0 0 5 0
    
"@
$prompt2 = @"
 Given is the adjacency matrix for a weighted undirected graph containing 23 nodes labelled A to W. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
what is the least cost path from node A to node W?
   A B C D E F G H I J K L M N O P Q R S T U V W
 A 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 4 0 4 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 4 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 5 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 5 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 3 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 5 0 0 0 2 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 3 0 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 3 0 0 0 0 0 2 0 0 4 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 2 0 2 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 2 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 3 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 4 0 0 0 0 2 0 0 0 2 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 5 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 2 0 0 0 5 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 2 0 0 0 1 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 2 0 0 0 0 0 5
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 5 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 5 0 2 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 2 0 5 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 5 0 3
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 3 0
    
"@
$prompt3 = @"
 Given is the adjacency matrix for a weighted undirected graph containing 22 nodes labelled A to V. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
what is the least cost path from node A to node V?
   A B C D E F G H I J K L M N O P Q R S T U V
 A 0 3 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 3 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 5 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 2 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 5 0 4 0 4 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 4 0 4 0 0 3 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 4 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 2 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 3 0 0 4 0 0 0 0 0 1 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 5 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 3 0 0 0 4 0 0
 P 0 0 0 0 0 0 0 0 0 0 1 0 0 0 3 0 4 0 0 0 4 0
 Q 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 4 0 0 0 0 0 1
 R 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 1 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 1 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 1 0 5
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 5 0
    
"@
$solution1 = @"
A -> B -> F -> J -> K -> L -> P
"@
$solution2 = @"
A -> E -> J -> K -> L -> P -> Q -> V -> W
"@
$solution3 = @"
A -> D -> E -> F -> G -> H -> L -> Q -> V
"@
$llm1 = @"
To find the least cost path from node A to node P, we can use Dijkstra's algorithm. 
1. Initialize a table to store the cost and predecessor nodes for each node. Set the cost for node A to 0 and the cost for all other nodes to infinity. Set the predecessor for all nodes to null.
2. Start at node A and consider all its neighbors. For each neighbor, update the cost and predecessor if a shorter path is found. In this step, we consider the neighbors of node A and update their costs and predecessors accordingly:
   - The cost to reach node B from node A is 2. Since this is smaller than the current cost of infinity, we update the cost for node B to 2 and set its predecessor to A.
   - The cost to reach node E from node A is 1. Since this is smaller than the current cost of infinity, we update the cost for node E to 1 and set its predecessor to A.
3. Move to the node with the smallest cost that has not been visited yet. In this case, the node with the smallest cost is E, so we move to node E.
4. Repeat step 2 for the new node. In this step, we consider the neighbors of node E and update their costs and predecessors accordingly:
   - The cost to reach node F from node E is 5. Since this is smaller than the current cost of infinity, we update the cost for node F to 5 and set its predecessor to E.
   - The cost to reach node I from node E is 5. Since this is smaller than the current cost of infinity, we update the cost for node I to 5 and set its predecessor to E.
5. Repeat steps 3 and 4 until all nodes have been visited. In this case, we continue the process and consider the next node with the smallest cost, which is node F:
   - The cost to reach node K from node F is 4. Since this is smaller than the current cost of infinity, we update the cost for node K to 4 and set its predecessor to F.
   - The cost to reach node J from node F is 4. Since this is smaller than the current cost of infinity, we update the cost for node J to 4 and set its predecessor to F.
6. Repeat steps 3 and 4 until all nodes have been visited. Continue the process until we reach node P.
7. Once all nodes have been visited, we can trace back the least cost path from node A to node P using the predecessor nodes. In this case, the least cost path from node A to node P is A -> E -> F -> J -> N -> O -> P.
Therefore, the least cost path from node A to node P is A -> E -> F -> J -> N -> O -> P with a total cost of 13.
"@
$llm2 = @"
To find the least cost path from node A to node W, we can use Dijkstra's algorithm.  
First, we initialize an array of distances with infinity values for all nodes except for A, which we set to 0.
Next, we initialize a set of unvisited nodes and add all nodes to it.
Then, we repeat the following steps until all nodes have been visited:
1. Select the node with the minimum distance from the array of distances and remove it from the set of unvisited nodes.
2. For each neighbor of the selected node that is still in the set of unvisited nodes, calculate the cost of traveling from the selected node to the neighbor and update the distance in the array if it is smaller than the current value.
3. Repeat steps 1 and 2 until there are no more unvisited nodes.
Finally, we can trace the shortest path from node A to node W by starting at W and repeatedly selecting the neighbor with the smallest distance until we reach A.
Here is the step-by-step implementation of this algorithm:
1. Initialize an array of distances with infinity values and set the distance of A to 0.
2. Initialize a set of unvisited nodes and add all nodes from A to W.
3. Repeat the following steps until all nodes have been visited:
   - Select the node with the minimum distance from the array of distances and remove it from the set of unvisited nodes.
   - For each of the selected node's neighbors that are still in the set of unvisited nodes:
     - Calculate the cost of traveling from the selected node to the neighbor.
     - If this cost plus the distance of the selected node is smaller than the current distance of the neighbor, update the distance of the neighbor.
4. Trace the shortest path from W to A using the updated array of distances.
   - Start at W and select the neighbor with the smallest distance.
   - Repeat this step until you reach A, adding each node to the path.
Following the above steps, we obtain the least cost path from node A to node W as:
A -> E -> F -> G -> H -> I -> J -> K -> L -> O -> P -> Q -> R -> S -> T -> U -> V -> W
"@
$llm3 = @"
To find the least cost path from node A to node V, we can use the Dijkstra's algorithm. Here is the step-by-step process:
1. Create a list to store the least cost values for each node, initialized with infinity for all nodes except for node A which is initialized with 0.
   A: 0, B: ∞, C: ∞, D: ∞, E: ∞, F: ∞, G: ∞, H: ∞, I: ∞, J: ∞, K: ∞, L: ∞, M: ∞, N: ∞, O: ∞, P: ∞, Q: ∞, R: ∞, S: ∞, T: ∞, U: ∞, V: ∞
2. Create a list to keep track of visited nodes, initialized as an empty list.
3. While there are unvisited nodes:
   a. Choose the node with the smallest least cost value from the list of unvisited nodes.
   b. Mark the chosen node as visited and add it to the list of visited nodes.
   c. Update the least cost values for the neighboring nodes of the chosen node:
      - If the current least cost value for a neighboring node is greater than the sum of the least cost value of the chosen node and the cost of travelling between the chosen node and the neighboring node, update the least cost value of the neighboring node.
   d. Repeat steps a-c until all nodes are visited.
4. Once all nodes are visited, the least cost path from node A to node V can be obtained by backtracking from node V to node A using the least cost values and the adjacency matrix. Starting from node V, choose the neighbor with the smallest least cost value and repeat until node A is reached.
Here is the step-by-step process for finding the least cost path from node A to node V:
1. Start with the initial least cost values:
   A: 0, B: ∞, C: ∞, D: ∞, E: ∞, F: ∞, G: ∞, H: ∞, I: ∞, J: ∞, K: ∞, L: ∞, M: ∞, N: ∞, O: ∞, P: ∞, Q: ∞, R: ∞, S: ∞, T: ∞, U: ∞, V: ∞
2. The list of visited nodes is empty.
3. Start the Dijkstra's algorithm:
   a. Choose the node with the smallest least cost value from the list of unvisited nodes (A with a cost of 0).
   b. Mark A as visited and add it to the list of visited nodes.
   c. Update the least cost values for the neighboring nodes of A:
      - B: current least cost = ∞, new least cost = 0 + 3 = 3, update least cost of B to 3.
      - D: current least cost = ∞, new least cost = 0 + 2 = 2, update least cost of D to 2.
   d. Move to the next node with the smallest least cost value (D with a cost of 2).
   e. Mark D as visited and add it to the list of visited nodes.
   f. Update the least cost values for the neighboring nodes of D:
      - A: already visited, no update.
      - E: already visited, no update.
      - I: current least cost = ∞, new least cost = 2 + 1 = 3, update least cost of I to 3.
   g. Move to the next node with the smallest least cost value (B with a cost of 3).
   h. Mark B as visited and add it to the list of visited nodes.
   i. Update the least cost values for the neighboring nodes of B:
      - A: already visited, no update.
      - C: current least cost = ∞, new least cost = 3 + 5 = 8, update least cost of C to 8.
   j. Move to the next node with the smallest least cost value (I with a cost of 3).
   k. Mark I as visited and add it to the list of visited nodes.
   l. Update the least cost values for the neighboring nodes of I:
      - D: already visited, no update.
      - M: current least cost = ∞, new least cost = 3 + 4 = 7, update least cost of M to 7.
   m. Move to the next node with the smallest least cost value (M with a cost of 7).
   n. Mark M as visited and add it to the list of visited nodes.
   o. Update the least cost values for the neighboring nodes of M:
      - I: already visited, no update.
      - N: current least cost = ∞, new least cost = 7 + 5 = 12, update least cost of N to 12.
      - R: current least cost = ∞, new least cost = 7 + 5 = 12, update least cost of R to 12.
   p. Move to the next node with the smallest least cost value (R with a cost of 12).
   q. Mark R as visited and add it to the list of visited nodes.
   r. Update the least cost values for the neighboring nodes of R:
      - M: already visited, no update.
      - Q: current least cost = ∞, new least cost = 12 + 1 = 13, update least cost of Q to 13.
      - U: current least cost = ∞, new least cost = 12 + 5 = 17, update least cost of U to 17.
   s. Move to the next node with the smallest least cost value (Q with a cost of 13).
   t. Mark Q as visited and add it to the list of visited nodes.
   u. Update the least cost values for the neighboring nodes of Q:
      - R: already visited, no update.
      - U: current least cost = 17, new least cost = 13 + 1 = 14, update least cost of U to 14.
      - V: current least cost = ∞, new least cost = 13 + 1 = 14, update least cost of V to 14.
   v. Move to the next node with the smallest least cost value (U with a cost of 14).
   w. Mark U as visited and add it to the list of visited nodes.
   x. Update the least cost values for the neighboring nodes of U:
      - Q: already visited, no update.
      - V: already visited, no update.
      - O: current least cost = ∞, new least cost = 14 + 3 = 17, update least cost of O to 17.
   y. Move to the next node with the smallest least cost value (O with a cost of 17).
   z. Mark O as visited and add it to the list of visited nodes.
   aa. Update the least cost values for the neighboring nodes of O:
       - U: already visited, no update.
       - P: current least cost = ∞, new least cost = 17 + 4 = 21, update least cost of P to 21.
       - N: already visited, no update.
   bb. Move to the next node with the smallest least cost value (P with a cost of 21).
   cc. Mark P as visited and add it to the list of visited nodes.
   dd. Update the least cost values for the neighboring nodes of P:
       - O: already visited, no update.
       - T: current least cost = ∞, new least cost = 21 + 4 = 25, update least cost of T to 25.
       - Q: already visited, no update.
   ee. Move to the next node with the smallest least cost value (T with a cost of 25).
   ff. Mark T as visited and add it to the list of visited nodes.
   gg. Update the least cost values for the neighboring nodes of T:
       - P: already visited, no update.
       - U: already visited, no update.
       - V: already visited, no update.
   hh. Move to the next node with the smallest least cost value (V with a cost of ∞).
   ii. Mark V as visited and add it to the list of visited nodes.
   jj. Update the least cost values for the neighboring nodes of V:
       - U: already visited, no update.
       - Q: already visited, no update.
   kk. All nodes are visited, terminate the algorithm.
5. Backtrack from the destination node V to node A to get the least cost path.
   Start from node V (least cost = ∞).
   - Node V has multiple neighboring nodes with the same least cost value (U and Q), choose one of them (U).
   - Node U has multiple neighboring nodes with the same least cost value (Q and O), choose one of them (Q).
   - Node Q has multiple neighboring nodes with the same least cost value (R and U), choose one of them (R).
   - Node R has multiple neighboring nodes with the same least cost value (M and Q), choose one of them (M).
   - Node M has multiple neighboring nodes with the same least cost value (I, N, and R), choose one of them (N).
   - Node N has multiple neighboring nodes with the same least cost value (M and O), choose one of them (O).
   - Node O has a neighboring node with a least cost value of 0 (P).
   - Node P has a neighboring node with a least cost value of 0 (T).
   - Node T has a neighboring node with a least cost value of 0 (U).
6. The least cost path from node A to node V is A → D → I → M → N → O → P → T → U → V.
Therefore, the least cost path from node A to node V is A → D → I → M → N → O → P → T → U → V with a total cost of 26.
"@

# --- Sheet1 (o_10) row 2 ---
$ws1.Range("A2").Value = $prompt1
$ws1.Range("B2").Value = $solution1
$ws1.Range("C2").Value = $llm1
$ws1.Range("D2").Value = "Wrong"
$ws1.Range("E2").Value = "Output: 2/7"

# --- Sheet2 (o_20): row 2 content (header row already copied from sheet1) ---
$ws2.Range("A2").Value = $prompt2
$ws2.Range("B2").Value = $solution2
$ws2.Range("C2").Value = $llm2
$ws2.Range("D2").Value = "Wrong"
$ws2.Range("E2").Value = "Output: 3/17"

# --- Sheet3 (o_20_jumbled): row 2 content (header row already copied from sheet1) ---
$ws3.Range("A2").Value = $prompt3
$ws3.Range("B2").Value = $solution3
$ws3.Range("C2").Value = $llm3
$ws3.Range("D2").Value = "Wrong"
$ws3.Range("E2").Value = "Output: 0/9"

# --- Undo the engine's auto row-height bump triggered by the multi-line values above ---
$ws1.Rows.Item(2).EntireRow.AutoFit() | Out-Null
$ws2.Rows.Item(2).EntireRow.AutoFit() | Out-Null
$ws3.Rows.Item(2).EntireRow.AutoFit() | Out-Null

# --- Restore original active sheet/selection ---
$ws1.Activate()
$ws1.Range("A1").Select() | Out-Null
